$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for data rows 2..516 all move from 2023-09-12 (45181)
# to 2023-09-13 (45182) - an increment of one day.
$ws.Range("C2:C516").Value2 = 45182
